$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" everywhere it appears ---
$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Narrow the (now shorter) Status-related columns ---
# Overview sheet: columns E ("zh-cn") and F ("de-de") mirror the Status column width
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de sheets: column C is "Status"
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
